$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "38.84", "1.00") that Excel would
# otherwise auto-convert to a real number, losing the original text formatting.
# Pre-format the Price column as Text, write the values, then restore the default
# (unstyled) look so the saved XML matches the original un-styled cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "45.972.34"
$ws.Cells.Item(2, 5).Value = "  -0.68%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.615.06"
$ws.Cells.Item(3, 5).Value = "  +0.05%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "309.58"
$ws.Cells.Item(5, 5).Value = "  -1.60%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "98.59"
$ws.Cells.Item(6, 5).Value = "  -2.21%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.75%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.04%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.579"
$ws.Cells.Item(9, 5).Value = "  -0.99%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "38.84"
$ws.Cells.Item(10, 5).Value = "  +0.08%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "OKB"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(11, 4).Value = "54.09"
$ws.Cells.Item(11, 5).Value = "  -0.43%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Dogecoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(12, 4).Value = "0.0840"
$ws.Cells.Item(12, 5).Value = "  -0.15%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "8.07"
$ws.Cells.Item(13, 5).Value = "  -3.13%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.011.25"
$ws.Cells.Item(14, 5).Value = "  +0.07%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +0.91%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.615.70"
$ws.Cells.Item(16, 5).Value = "  +0.26%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.03%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "14.80"
$ws.Cells.Item(18, 5).Value = "  -2.54%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "45.987.48"
$ws.Cells.Item(19, 5).Value = "  -1.06%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -1.05%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "6.78"
$ws.Cells.Item(21, 5).Value = "  +0.41%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "12.78"
$ws.Cells.Item(22, 5).Value = "  -4.30%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "74.36"
$ws.Cells.Item(23, 5).Value = "  +4.91%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "284.29"
$ws.Cells.Item(24, 5).Value = "  +11.49%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -1.79%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "2.27"
$ws.Cells.Item(26, 5).Value = "  +2.67%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +5.24%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "0.999"
$ws.Cells.Item(28, 5).Value = "  -0.03%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.31%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "10.55"
$ws.Cells.Item(30, 5).Value = "  -0.88%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "38.81"
$ws.Cells.Item(31, 5).Value = "  -4.76%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -2.87%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.06%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "3.61"
$ws.Cells.Item(34, 5).Value = "  -3.51%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.52%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "157.05"
$ws.Cells.Item(36, 5).Value = "  +2.03%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "2.81"
$ws.Cells.Item(37, 5).Value = "  -2.34%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.74%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "0.121"
$ws.Cells.Item(39, 5).Value = "  +1.92%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.25%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "15.94"
$ws.Cells.Item(41, 5).Value = "  -7.25%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "0.0327"
$ws.Cells.Item(42, 5).Value = "  -0.28%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "EnergySwap"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(43, 4).Value = "21.57"
$ws.Cells.Item(43, 5).Value = "  +2.86%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "NEARProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(44, 4).Value = "3.54"
$ws.Cells.Item(44, 5).Value = "  -2.63%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "RenderToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(45, 4).Value = "4.03"
$ws.Cells.Item(45, 5).Value = "  -4.98%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "2.106.90"
$ws.Cells.Item(46, 5).Value = "  +3.53%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "1.00"
$ws.Cells.Item(47, 5).Value = "  +0.11%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "93.68"
$ws.Cells.Item(48, 5).Value = "  +2.81%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "9.19"
$ws.Cells.Item(49, 5).Value = "  -0.55%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "109.28"
$ws.Cells.Item(50, 5).Value = "  -2.88%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -1.09%  "

# Restore the Price column to the workbook-default (unstyled) appearance now
# that the values are safely stored as text.
$priceRange.Style = "Normal"
